# The commit swaps the deck's theme ("Integral" / Red-Violet colour
# scheme) for the plain default "Office Theme" colour scheme — i.e. the
# author picked a different design from the Design gallery.
#
# ppt/theme/theme1.xml is the theme actually driving every slide (it's
# the part referenced by the single slide master), so that's what needs
# its 12 colour-scheme slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) rewritten to the target "Office" values. The font scheme and
# format scheme (fills/lines/effects) are already byte-identical between
# the two themes in this deck, so only the colours need to change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexToComRgb($officeThemeColors[$i])
}
